$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 796
$ws.Range("F6").Value = 250
$ws.Range("F15").Value = 714
$ws.Range("F16").Value = 841
$ws.Range("F17").Value = 9107
$ws.Range("F27").Value = 271
$ws.Range("F28").Value = 471
$ws.Range("F33").Value = 49
$ws.Range("F37").Value = 163

# Sheet "演出" (performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 67

# Sheet "本地生活" (local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 799

# Sheet "全部类型" (all types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 799
$ws.Range("F6").Value = 796
$ws.Range("F10").Value = 250
$ws.Range("F22").Value = 841
$ws.Range("F23").Value = 9110
$ws.Range("F32").Value = 271
$ws.Range("F33").Value = 471
$ws.Range("F35").Value = 67
$ws.Range("F36").Value = 67
$ws.Range("F42").Value = 49
